# "fixed tested points & test save location"
#
# Refresh the 10 "tested point" rows: columns A/B/C get corrected
# measurement values, every row is made a touch taller, and the font
# used throughout the data grid is pinned to an explicit black (instead
# of following the workbook's "theme 1 / text 1" automatic color) so the
# saved file renders identically regardless of theme.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tested-point values for columns A / B / C, one row per entry
# (row, A, B, C).
$rows = @(
    @(1,  0.04,  -0.13, -0.34),
    @(2,  -0.27, -0.05,  0.3),
    @(3,  0.07,  -0.24, -0.27),
    @(4,  0.27,  -0.11,  0.17),
    @(5,  -0.03, -0.11,  0.77),
    @(6,  0.06,   0.03,  0.16),
    @(7,  0.17,  -0.09, -0.18),
    @(8,  0.22,  -0.44,  0.46),
    @(9,  -0.35, -0.33,  0.39),
    @(10, 0.08,  -0.03,  0.29)
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]

    # Rows grew from 16.5pt to 18.75pt.
    $ws.Rows.Item($r).RowHeight = 18.75
}

# Pin the data font to explicit black (RGB 0,0,0) instead of the
# "theme 1" automatic text color.
$ws.Range("A1:F10").Font.Color = 0
